$d = $word.ActiveDocument

# The original document ends with a paragraph "...screen" that carries the
# hidden "_GoBack" bookmark right at the very end of the document. Word
# relocates that bookmark to the last edited position as you type, so once
# we are done adding the new log entries it needs to end up immediately
# after "...animations " in the final new paragraph.
#
# The COM surface here has an edge case where Bookmarks.Add / Range
# operations that land exactly on the document's absolute end position
# misbehave, so rather than fight that boundary we remove the old bookmark
# up front and re-introduce bookmarkStart/bookmarkEnd markers as literal
# OOXML (via InsertXML) in the correct spot alongside the rest of the new
# content.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Turn "... The bullets now appear correctly on the screen" into its own
# trailing "." run (mirrors the diff's new <w:r><w:t>.</w:t></w:r>).
$lastPara = $d.Paragraphs.Last
$endPos = $lastPara.Range.End
$dotRange = $d.Range($endPos, $endPos)
$dotRange.InsertAfter(".")

$newParasXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Added gravity </w:t></w:r><w:r><w:t xml:space="preserve">so the player will always be walking on the </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>floor</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> so you don’t need to position them on it manually</w:t></w:r><w:r><w:t>, it’s just a set negative number as realistically there isn’t any jumping or falling in this kind of game anyway and it really only there to make sure the player starts on the ground</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Fixed small error with touch manager where it was trying to grab components off a null object if you didn’t tap a game object (wasn’t causing any issues with movement, was just putting out an error that went unnoticed)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Fixed the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>thumbsticks</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> not </w:t></w:r><w:r><w:t>being anchored to the corners of the screen</w:t></w:r><w:r><w:t>/made them bigger</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Added a laser sound to the bullets, made it play at a random pitch between 0.75 and 1.</w:t></w:r><w:r><w:t>5 so it wasn’t repetitive to hear</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Swapped back to using the sci fi </w:t></w:r><w:r><w:t>character as I wanted to test out more animation layer masks, using the soldiers run anim</w:t></w:r><w:r><w:t>ation</w:t></w:r><w:r><w:t xml:space="preserve"> on a separate layer and the rest are animations that came with the sci fi character. This lets me use the much nicer run that otherwise would make the sci fi characters arms fold in on themselves. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Next</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> I’m going to set up the level to look a lot nicer with some free unity assets.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Set up the level with some free sci fi assets, am now going to put in some robots that will stand around for you to shoot. So </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>far</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> they </w:t></w:r><w:r><w:t xml:space="preserve">have an idle animation that transitions to a hit </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>one when you shoot them, am going to implement some code that will cause them to die when shot too many times.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Enemy robots now take one damage per bullet and will have their collider deactivated as well as play the death animation when their hp hits zero.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Fixed an issue where the enemies were being knocked back by bullets</w:t></w:r><w:r><w:t xml:space="preserve"> by turning the root motion off their animations </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

# Insert every new changelog bullet (with blank ListParagraph spacer
# paragraphs in between, proofErr spell/grammar markers, the mid-sentence
# page-break marker, and the relocated _GoBack bookmark on the very last
# run) right after the paragraph we just appended "." to.
$lastPara2 = $d.Paragraphs.Last
$insertPos = $lastPara2.Range.End
$insertRange = $d.Range($insertPos, $insertPos)
[void]$insertRange.InsertXML($newParasXml)
